# Update the "想去人数" (want-to-go count) figures in column F for both the
# "展览" sheet and the "全部类型" sheet, matching the refreshed data pull.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1094
    5  = 4687
    7  = 403
    8  = 1410
    11 = 1215
    13 = 674
    15 = 59
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
